$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand-new product row right before the "ZINCTRON 30 CAPS" row
# (worksheet row 81), pushing it and everything after it down by one row.
# ---------------------------------------------------------------------------
$ws.Rows("81:81").Insert($excel.XlInsertShiftDirection.xlShiftDown)

# ---------------------------------------------------------------------------
# Re-create the per-column formatting for the new row 81 (the plain insert
# above only shifts existing rows - which keep their formatting/merges - but
# the brand-new blank row has none of that, so it is rebuilt by hand here to
# match its neighbours).
# ---------------------------------------------------------------------------
$ws.Rows.Item(81).RowHeight = 25.5

# Column A:B  (index number) -------------------------------------------------
$rngA = $ws.Range("A81:B81")
$rngA.Font.Name = "Mega"
$rngA.Font.Size = 9
$rngA.Font.Color = 16777215
$rngA.Interior.Pattern = $excel.XlPattern.xlPatternSolid
$rngA.Interior.Color = 33023
$rngA.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
$rngA.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
$rngA.WrapText = $true
$rngA.Borders.Item(9).LineStyle = 1
$rngA.Borders.Item(9).Weight = 2
$rngA.Borders.Item(9).Color = 13882323

# Columns C:G, N:O (item name / price text) ----------------------------------
$rngCG = $ws.Range("C81:G81")
$rngN  = $ws.Range("N81:O81")
foreach ($rng in @($rngCG, $rngN)) {
    $rng.NumberFormat = "@"
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 9
    $rng.Font.Color = 0
    $rng.Interior.Pattern = $excel.XlPattern.xlPatternNone
    $rng.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
    $rng.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
    $rng.WrapText = $true
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = 13882323
}

# Columns H:K (current balance) ----------------------------------------------
$rngH = $ws.Range("H81:K81")
$rngH.NumberFormat = "@"
$rngH.Font.Name = "Arial"
$rngH.Font.Size = 9
$rngH.Font.Color = 0
$rngH.Interior.Pattern = $excel.XlPattern.xlPatternSolid
$rngH.Interior.Color = 16119285
$rngH.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
$rngH.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
$rngH.WrapText = $true
$rngH.Borders.Item(9).LineStyle = 1
$rngH.Borders.Item(9).Weight = 2
$rngH.Borders.Item(9).Color = 13882323

# Columns L:M (order limit) ---------------------------------------------------
$rngL = $ws.Range("L81:M81")
$rngL.NumberFormat = '#,##0.##;"["#,##0.##"]";0'
$rngL.Font.Name = "Arial"
$rngL.Font.Size = 9
$rngL.Font.Color = 0
$rngL.Interior.Pattern = $excel.XlPattern.xlPatternNone
$rngL.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
$rngL.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
$rngL.WrapText = $true
$rngL.Borders.Item(9).LineStyle = 1
$rngL.Borders.Item(9).Weight = 2
$rngL.Borders.Item(9).Color = 13882323

# Column P (sell price) -------------------------------------------------------
$rngP = $ws.Range("P81")
$rngP.NumberFormat = "0.00"
$rngP.Font.Name = "Arial"
$rngP.Font.Size = 9
$rngP.Font.Color = 0
$rngP.Interior.Pattern = $excel.XlPattern.xlPatternNone
$rngP.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
$rngP.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
$rngP.WrapText = $true
$rngP.Borders.Item(9).LineStyle = 1
$rngP.Borders.Item(9).Weight = 2
$rngP.Borders.Item(9).Color = 13882323

# Column Q (number of transactions) ------------------------------------------
$rngQ = $ws.Range("Q81")
$rngQ.NumberFormat = "@"
$rngQ.Font.Name = "Arial"
$rngQ.Font.Size = 9
$rngQ.Font.Color = 0
$rngQ.Interior.Pattern = $excel.XlPattern.xlPatternNone
$rngQ.HorizontalAlignment = $excel.XlHAlign.xlHAlignCenter
$rngQ.VerticalAlignment = $excel.XlVAlign.xlVAlignCenter
$rngQ.WrapText = $true
$rngQ.Borders.Item(9).LineStyle = 1
$rngQ.Borders.Item(9).Weight = 2
$rngQ.Borders.Item(9).Color = 13882323

# Re-create the merges for the new row (same layout as every other data row)
$ws.Range("A81:B81").Merge()
$ws.Range("C81:G81").Merge()
$ws.Range("H81:K81").Merge()
$ws.Range("L81:M81").Merge()
$ws.Range("N81:O81").Merge()

# ---------------------------------------------------------------------------
# Populate the new row with the new product's data.
# ---------------------------------------------------------------------------
$ws.Range("A81").Value = 75
$ws.Range("C81").Value = "ZESTRIL 20MG 10 TAB"
$ws.Range("H81").Value = "2:0"
$ws.Range("L81").Value = "1"
$ws.Range("N81").Value = "68.00"
$ws.Range("P81").Value = "68.0000"
$ws.Range("Q81").Value = "1:0"

# ---------------------------------------------------------------------------
# Update the running total (now on row 93) to include the new item's price.
# ---------------------------------------------------------------------------
$ws.Range("P93").Value = 4612.01

# ---------------------------------------------------------------------------
# Refresh the "printed at" timestamp shown in the footer (now row 94).
# ---------------------------------------------------------------------------
$ws.Range("G94").Value = "Thursday, 7 August, 2025 7:51 PM"
